$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "MAIN Config." -- default code config.sh files, paths now relative
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Write order chosen to match the original author's editing sequence (keeps
# the shared-string table append order identical to the source edit).
$ws1.Range("B7").Value = "openmc_config.sh"
$ws1.Range("B3").Value = "mcnp_config.sh"
$ws1.Range("B5").Value = "serpent_config.sh"
$ws1.Range("B4").Value = "/home/sbradnam/Software/freia/Serpent2_src/v2.1.32_ccfe/sss2"
$ws1.Range("B6").Value = "/home/sbradnam/Software/freia/OPENMC_311022/openmc/build/bin/openmc"

# B4 (new Serpent executable row) picks up the same style as the other
# "extra" rows (B8/B9) further down the sheet.
$ws1.Range("B8").Copy()
$ws1.Range("B4").PasteSpecial(-4122)

# The placeholder cell for "d1S executable" (row 8) is removed entirely.
$ws1.Range("B8").Clear()

$ws1.Range("B6").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Computational benchmarks" -- Sphere Leakage Test now also runs
# under Serpent (and OnlyInput flips to true)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# D4 already holds the text "true" -- copy its value (only) onto C4/E4 so the
# result stays a shared-string "true" rather than a native Boolean TRUE.
$ws2.Range("D4").Copy()
$ws2.Range("C4").PasteSpecial(-4163)  # xlPasteValues
$ws2.Range("D4").Copy()
$ws2.Range("E4").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$ws2.Range("D4").Select()

# ---------------------------------------------------------------------------
# Sheet 4: "Libraries" -- add a dedicated Serpent xsdir column (mirrors the
# MCNP column for now)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("E1").Value = "Serpent"
$ws4.Range("E2").Value = "/home/mcnp/xs/xsdir_mcnp6.2"
$ws4.Range("E3").Value = "/home/mcnp/xs/xsdir_mcnp6.2_old"
$ws4.Range("E4").Value = "/home/mcnp/xs/xsdir_mcnp6.2_jeff33_endfb71_fendl32b_irdff105_tt"
$ws4.Range("E5").Value = "/home/mcnp/xs/xsdir_mcnp6.2_fendl32b_rw"
$ws4.Range("E6").Value = "/home/mcnp/xs/xsdir_mcnp6.2"
$ws4.Range("E7").Value = "/home/mcnp/xs/xsdir_mcnp6.2_endfb8"

$ws4.Columns.Item(5).AutoFit()

# Libraries becomes the active tab/selection.
$ws4.Range("E2").Select()
$ws4.Activate()

Write-Output "edit complete"
